$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.370.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.261.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '119.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.639'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.73%  '
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0937'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.899'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.604.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.266.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.388.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("E22").Value = '  -5.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.56%  '
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0910'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.82%  '
$ws.Range("E36").Value = '  +1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0377'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.240'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.95%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '75.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +35.83%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.11%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.658'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +16.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.90%  '
